# Updated symbol list refresh (Mon Feb 13 19:11:21 UTC 2023).
# Source cells in this sheet are plain text (inlineStr) even though many
# values look numeric (prices, percent changes, the "Hora" hour column).
# To keep them as text through COM - rather than have Excel auto-convert
# a numeric-looking string to a real number - each cell is briefly
# switched to the "Text" number format before the value is written, then
# restored to the default "Normal" style so no extra formatting lingers
# on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @"
D2	288.06
E2	-10.06%
G2	19
D3	40.35
E3	-2.40%
G3	19
D4	5.039
E4	-4.18%
G4	19
E5	-5.67%
G5	19
D6	4.283
E6	-1.37%
G6	19
D7	1.511
E7	-12.75%
G7	19
D8	0.9191
E8	-2.70%
G8	19
E9	-4.82%
G9	19
D10	0.1726
E10	-6.82%
G10	19
D11	0.08633
E11	-5.74%
G11	19
D12	0.04170
E12	-3.26%
G12	19
D13	0.1053
E13	0.16%
G13	19
E14	-1.89%
G14	19
D15	0.005913
E15	-1.14%
G15	19
D16	3.396
E16	1.59%
G16	19
G17	19
D18	0.3288
E18	-1.36%
G18	19
D19	7.827
E19	0.84%
G19	19
E20	-0.30%
G20	19
E21	1.93%
G21	19
D22	0.03860
E22	-4.12%
G22	19
E23	-0.16%
G23	19
D24	0.003840
E24	-6.59%
G24	19
D25	0.0001281
E25	0.76%
G25	19
D26	0.0003721
G26	19
G27	19
G28	19
G29	19
G30	19
G31	19
G32	19
G33	19
G34	19
G35	19
G36	19
G37	19
D38	0.02316
E38	-9.29%
G38	19
D39	0.04980
E39	-6.86%
G39	19
D40	0.006553
E40	228.85%
G40	19
D41	0.007699
E41	-0.73%
G41	19
D42	0.1275
G42	19
D43	0.007349
E43	-0.23%
G43	19
E44	-14.74%
G44	19
D45	0.3119
E45	-2.04%
G45	19
D46	0.00006420
E46	-4.72%
G46	19
E47	-0.17%
G47	19
D48	0.2424
E48	20.86%
G48	19
E49	-0.10%
G49	19
D50	0.00002100
E50	-0.17%
G50	19
D51	0.0002000
E51	-0.17%
G51	19
"@

$rows = $updates -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
foreach ($line in $rows) {
    $parts = $line -split "`t"
    $cellRef = $parts[0].Trim()
    $newValue = $parts[1]

    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Write-Output "Updated $($rows.Count) cells (Price/Volume(1h)/Hora columns) with refreshed crypto data."
